$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 200657
$ws.Range("N3").Value = -200885
$ws.Range("J3").Value = 200657
$ws.Range("L3").Value = 200657
$ws.Range("N19").Value = -1126.5
$ws.Range("J19").Value = 776.5
$ws.Range("L19").Value = 776.5
$ws.Range("H19").Value = 809.8
$ws.Range("N64").Value = -5496
$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("H67").Value = 5000
$ws.Range("N67").Value = -6716
$ws.Range("L67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("N88").Value = -67276.336
$ws.Range("K88").Value = 363.25
$ws.Range("H88").Value = 28692.285
$ws.Range("J88").Value = 66464.336
$ws.Range("I88").Value = 363.25
$ws.Range("L88").Value = 66464.336
$ws.Range("M88").Value = 42.75
$ws.Range("I91").Value = 363.25
$ws.Range("J91").Value = 66464.336
$ws.Range("H91").Value = 28692.285
$ws.Range("M91").Value = 1040.75
$ws.Range("K91").Value = 363.25
$ws.Range("L91").Value = 66464.336
$ws.Range("N91").Value = -69272.336
$ws.Range("N95").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("H102").Value = 200657
$ws.Range("J102").Value = 200657
$ws.Range("L102").Value = 200657
$ws.Range("N102").Value = -207147
$ws.Range("L124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("H124").Value = 75000
$ws.Range("N124").Value = -84820
$ws.Range("N128").Value = -69960
$ws.Range("H128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("I132").Value = 2010.1154
$ws.Range("K132").Value = 6030.3462
$ws.Range("M132").Value = -3500.3462
$ws.Range("H132").Value = 2129
$ws.Range("H136").Value = 67065
$ws.Range("H137").Value = 2136.9644
$ws.Range("I137").Value = 1675.9474
$ws.Range("M137").Value = -2477.8422
$ws.Range("K137").Value = 5027.8422

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 577.625
$ws.Range("K2").Value = 511.2
$ws.Range("M2").Value = -398.2
$ws.Range("I2").Value = 511.2
$ws.Range("H45").Value = 1750
$ws.Range("M45").Value = -623
$ws.Range("K45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("H102").Value = 1985.7142
$ws.Range("I102").Value = 1985.7142
$ws.Range("K102").Value = 1985.7142
$ws.Range("M102").Value = -363.7141999999999
$ws.Range("H116").Value = 577.625
$ws.Range("M116").Value = 1782.8
$ws.Range("K116").Value = 511.2
$ws.Range("I116").Value = 511.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 511.2
$ws.Range("K3").Value = 511.2
$ws.Range("H3").Value = 577.625
$ws.Range("M3").Value = -397.2
$ws.Range("J94").Value = 1210
$ws.Range("L94").Value = 1210
$ws.Range("I94").Value = 1377
$ws.Range("K94").Value = 1377
$ws.Range("N94").Value = -2112
$ws.Range("H94").Value = 1343.6
$ws.Range("M94").Value = -926
$ws.Range("I99").Value = 3569.6365
$ws.Range("H99").Value = 3597.0833
$ws.Range("K99").Value = 3569.6365
$ws.Range("M99").Value = -2071.6365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J31").Value = 2704.3333
$ws.Range("H31").Value = 1965.4166
$ws.Range("M31").Value = -1424.1111
$ws.Range("L31").Value = 2704.3333
$ws.Range("I31").Value = 1719.1111
$ws.Range("N31").Value = -3294.3333
$ws.Range("K31").Value = 1719.1111
$ws.Range("M34").Value = -1517.1111
$ws.Range("N34").Value = -3108.3333
$ws.Range("I34").Value = 1719.1111
$ws.Range("K34").Value = 1719.1111
$ws.Range("J34").Value = 2704.3333
$ws.Range("L34").Value = 2704.3333
$ws.Range("H34").Value = 1965.4166
$ws.Range("I58").Value = 1659.5
$ws.Range("K58").Value = 1659.5
$ws.Range("H58").Value = 3677.6667
$ws.Range("M58").Value = -1456.5
$ws.Range("I107").Value = 406.08334
$ws.Range("K107").Value = 406.08334
$ws.Range("M107").Value = 1513.91666
$ws.Range("H107").Value = 406.08334
$ws.Range("L132").Value = 14925
$ws.Range("N132").Value = -19985
$ws.Range("J132").Value = 4975
$ws.Range("H132").Value = 3524
$ws.Range("H136").Value = 3677.6667
$ws.Range("K136").Value = 4978.5
$ws.Range("M136").Value = -2428.5
$ws.Range("I136").Value = 1659.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N92").Value = -5941.9998
$ws.Range("H92").Value = 1061.5
$ws.Range("J92").Value = 1148.6666
$ws.Range("L92").Value = 3445.9998
$ws.Range("I107").Value = 419.6
$ws.Range("N107").Value = -6438
$ws.Range("K107").Value = 1258.8
$ws.Range("M107").Value = 661.1999999999998
$ws.Range("J107").Value = 866
$ws.Range("L107").Value = 2598
$ws.Range("H107").Value = 587
$ws.Range("J114").Value = 0
$ws.Range("M114").Value = 2654
$ws.Range("K114").Value = 600
$ws.Range("L114").Value = 0
$ws.Range("I114").Value = 200
$ws.Range("N114").ClearContents()
$ws.Range("H114").Value = 200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("I36").Value = 25000
$ws.Range("H36").Value = 25000
$ws.Range("K36").Value = 25000
$ws.Range("M36").Value = -24515
$ws.Range("N40").ClearContents()
$ws.Range("J40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N43").Value = -10302
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10000
$ws.Range("M43").ClearContents()
$ws.Range("L43").Value = 10000
$ws.Range("K43").Value = 0
$ws.Range("H80").Value = 1636.2727
$ws.Range("I80").Value = 1400
$ws.Range("K80").Value = 1400
$ws.Range("M80").Value = -402
$ws.Range("H83").Value = 1636.2727
$ws.Range("K83").Value = 7000
$ws.Range("I83").Value = 1400
$ws.Range("M83").Value = -2008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I93").Value = 3166.1667
$ws.Range("L93").Value = 2997
$ws.Range("H93").Value = 3142
$ws.Range("J93").Value = 2997
$ws.Range("K93").Value = 3166.1667
$ws.Range("N93").Value = -5493
$ws.Range("M93").Value = -1918.1667
$ws.Range("I132").Value = 15773
$ws.Range("K132").Value = 47319
$ws.Range("M132").Value = -44789
$ws.Range("H132").Value = 15775
$ws.Range("H136").Value = 3197.3333
$ws.Range("K136").Value = 7310.400000000001
$ws.Range("M136").Value = -4760.400000000001
$ws.Range("I136").Value = 2436.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2580.8
$ws.Range("J62").Value = 3221
$ws.Range("L62").Value = 3221
$ws.Range("N62").Value = -4469
$ws.Range("L65").Value = 16105
$ws.Range("H65").Value = 2580.8
$ws.Range("N65").Value = -22345
$ws.Range("J65").Value = 3221
$ws.Range("J75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("L75").Value = 25000
$ws.Range("H75").Value = 25000
$ws.Range("N78").Value = -84360
$ws.Range("H78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("J78").Value = 25000
$ws.Range("K122").Value = 7033.1535
$ws.Range("M122").Value = -4583.1535
$ws.Range("J122").Value = 2350.5
$ws.Range("H122").Value = 2345.2
$ws.Range("N122").Value = -11951.5
$ws.Range("I122").Value = 2344.3845
$ws.Range("L122").Value = 7051.5
